$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row is inserted at row 126, pushing all subsequent
# records (old rows 126-239) down by one row (new rows 127-240).
$ws.Rows("126:126").Insert()

# Fill in the newly inserted row 126 with this week's record.
$ws.Cells.Item(126, 1).Value = 5
$ws.Cells.Item(126, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(126, 3).Value = "Maule"
$ws.Cells.Item(126, 4).Value = 44484
$ws.Cells.Item(126, 5).Value = 7
$ws.Cells.Item(126, 6).Value = 100112043
$ws.Cells.Item(126, 7).Value = "Pepino ensalada"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 300
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 13).Value = 12000
$ws.Cells.Item(126, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(126, 15).Value = "Región del Maule"
$ws.Cells.Item(126, 16).Value = 150
$ws.Cells.Item(126, 17).Value = 80
$ws.Cells.Item(126, 18).Value = "Hortaliza"
